$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 61, pushing existing rows 61-163 down to 62-164
$ws.Rows("61:61").Insert()

# Populate the newly inserted row 61 with the new record's data.
# Static/template columns match the rest of the dataset.
$ws.Range("A61").Value = 5
$ws.Range("B61").Value = "Macroferia Regional de Talca"
$ws.Range("C61").Value = "Maule"
$ws.Range("D61").Value = 44725
$ws.Range("D61").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E61").Value = 7
$ws.Range("F61").Value = 100112031
$ws.Range("G61").Value = "Poroto verde"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 150
$ws.Range("K61").Value = 25000
$ws.Range("L61").Value = 25000
$ws.Range("M61").Value = 25000
$ws.Range("N61").Value = "`$/malla 25 kilos"
$ws.Range("O61").Value = "Región de Arica y Parinacota"
$ws.Range("P61").Value = 1000
$ws.Range("Q61").Value = 25
$ws.Range("R61").Value = "Hortaliza"
